# Applies the "gh-pages output" refresh: updates the "想去人数" (column F)
# counts on the 展览 / 演出 / 本地生活 / 全部类型 sheets to their newer values.

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F15").Value = 1761
$ws.Range("F16").Value = 41
$ws.Range("F18").Value = 32
$ws.Range("F19").Value = 51
$ws.Range("F20").Value = 656
$ws.Range("F21").Value = 405
$ws.Range("F22").Value = 720
$ws.Range("F23").Value = 78958
$ws.Range("F24").Value = 78958
$ws.Range("F27").Value = 33619
$ws.Range("F28").Value = 33619
$ws.Range("F29").Value = 501
$ws.Range("F34").Value = 960
$ws.Range("F35").Value = 288
$ws.Range("F37").Value = 579
$ws.Range("F38").Value = 972
$ws.Range("F42").Value = 446
$ws.Range("F46").Value = 382

# --- 演出 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 1730
$ws.Range("F16").Value = 19
$ws.Range("F18").Value = 75
$ws.Range("F21").Value = 69
$ws.Range("F22").Value = 73

# --- 本地生活 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 85

# --- 全部类型 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F19").Value = 1761
$ws.Range("F21").Value = 75
$ws.Range("F22").Value = 32
$ws.Range("F23").Value = 51
$ws.Range("F24").Value = 656
$ws.Range("F26").Value = 405
$ws.Range("F27").Value = 720
$ws.Range("F28").Value = 78958
$ws.Range("F30").Value = 33619
$ws.Range("F31").Value = 501
$ws.Range("F36").Value = 960
$ws.Range("F38").Value = 288
$ws.Range("F39").Value = 579
$ws.Range("F40").Value = 972
$ws.Range("F41").Value = 972
$ws.Range("F47").Value = 446
$ws.Range("F49").Value = 382

$wb.Save()
